# Started normalization_rna module. Added rasterisation.
# Adds a new "pbmc_3k" / SmartSeq2 RNA dataset row to the 10x PBMC datasets
# sample sheet, and updates the current selection accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data row (row 4): sample, experiment, technology, assays, path, ...
$ws.Range("A4").Value = "pbmc"
$ws.Range("B4").Value = "pbmc_3k"
$ws.Range("C4").Value = "smartseq2"
$ws.Range("D4").Value = "RNA"
$ws.Range("E4").Value = "datasets/10x_SmartSeq2_pbmc_GSE132044/counts/smartseq2/counts_table.tsv.gz"
$ws.Range("I4").Value = "ENSEMBL"

# Move/extend the active selection the way it ended up after the edit.
$ws.Range("D8").Select() | Out-Null
